$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize the "Source" column so that all the paleo-reconstruction rows
# that previously cited later Meko vintages now cite "Meko et al 2017"
# (the later "Meko et al 2018/2019/2020" sources are no longer used and
# will drop out of the shared-strings table on save).
$ws.Range("E6").Value = "Meko et al 2017"
$ws.Range("E11").Value = "Meko et al 2017"
$ws.Range("E12").Value = "Meko et al 2017"
$ws.Range("E13").Value = "Meko et al 2017"
$ws.Range("E18").Value = "Meko et al 2017"

# Scroll the view down and leave the selection on E18, matching where the
# author ended up working in the sheet.
$ws.Range("E18").Select()
$excel.ActiveWindow.ScrollRow = 10
